# Update "想去人数" (interest count) and "最低票价" (min price) figures
# to the latest scraped values, for sheets "展览" and "全部类型".
# (Sheets "演出" and "本地生活" are unaffected by this data refresh.)

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 482
$ws1.Range("F3").Value = 5741
$ws1.Range("G3").Value = 62
$ws1.Range("F5").Value = 73
$ws1.Range("F6").Value = 97
$ws1.Range("F7").Value = 4
$ws1.Range("F8").Value = 56
$ws1.Range("F9").Value = 544

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 482
$ws4.Range("F3").Value = 5741
$ws4.Range("G3").Value = 62
$ws4.Range("F6").Value = 73
$ws4.Range("F7").Value = 97
$ws4.Range("F8").Value = 4
$ws4.Range("F10").Value = 56
$ws4.Range("F11").Value = 544
